$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.5837342739105225
$ws.Cells.Item(2, 5).Value = 4074.437543259144
$ws.Cells.Item(2, 6).Value = 0.115228378485947
$ws.Cells.Item(2, 7).Value = 0.09427896080783064
$ws.Cells.Item(2, 8).Value = 0.09164862871090457
$ws.Cells.Item(2, 9).Value = 0.09164862871090457
$ws.Cells.Item(2, 10).Value = 0.09162872496323575
$ws.Cells.Item(2, 11).Value = 0.09162872496323575
$ws.Cells.Item(2, 12).Value = 0.09101188561319042
$ws.Cells.Item(2, 13).Value = 0.09101188561319042
$ws.Cells.Item(2, 14).Value = 0.09073035137147749
$ws.Cells.Item(2, 15).Value = 0.09073035137147749
$ws.Cells.Item(2, 16).Value = 0.09071455732584895
$ws.Cells.Item(2, 17).Value = 0.09071455732584895
$ws.Cells.Item(2, 18).Value = 0.09071455732584895
$ws.Cells.Item(2, 19).Value = 0.09071455732584895
$ws.Cells.Item(2, 20).Value = 0.09071455732584895
$ws.Cells.Item(2, 21).Value = 0.09053651821903502
$ws.Cells.Item(2, 22).Value = 0.09022796526948479
$ws.Cells.Item(2, 23).Value = 0.09022796526948479
$ws.Cells.Item(2, 24).Value = 0.09015524345774807
$ws.Cells.Item(2, 25).Value = 0.09009040045339461
$ws.Cells.Item(3, 3).Value = 0.5937159061431885
$ws.Cells.Item(3, 5).Value = 4010.74386216471
$ws.Cells.Item(3, 6).Value = 0.1105902634219777
$ws.Cells.Item(3, 7).Value = 0.09189161393558125
$ws.Cells.Item(3, 8).Value = 0.09176008186304557
$ws.Cells.Item(3, 9).Value = 0.09176008186304557
$ws.Cells.Item(3, 10).Value = 0.0907168224867408
$ws.Cells.Item(3, 11).Value = 0.08975750645070747
$ws.Cells.Item(3, 12).Value = 0.08975750645070747
$ws.Cells.Item(3, 13).Value = 0.08975750645070747
$ws.Cells.Item(3, 14).Value = 0.08975750645070747
$ws.Cells.Item(3, 15).Value = 0.08975750645070747
$ws.Cells.Item(3, 16).Value = 0.08962920808997456
$ws.Cells.Item(3, 17).Value = 0.08962920808997456
$ws.Cells.Item(3, 18).Value = 0.08962920808997456
$ws.Cells.Item(3, 19).Value = 0.08962920808997456
$ws.Cells.Item(3, 20).Value = 0.08938708446695913
$ws.Cells.Item(3, 21).Value = 0.08916406507233701
$ws.Cells.Item(3, 22).Value = 0.08916406507233701
$ws.Cells.Item(3, 23).Value = 0.08909089488286877
$ws.Cells.Item(3, 24).Value = 0.08898035547547972
$ws.Cells.Item(3, 25).Value = 0.08884880822933156
$ws.Cells.Item(4, 3).Value = 0.6720044612884521
$ws.Cells.Item(4, 5).Value = 3982.49422706413
$ws.Cells.Item(4, 6).Value = 0.1117251618707877
$ws.Cells.Item(4, 7).Value = 0.09371622492156641
$ws.Cells.Item(4, 8).Value = 0.09188442673609226
$ws.Cells.Item(4, 9).Value = 0.08868758830594635
$ws.Cells.Item(4, 10).Value = 0.08868758830594635
$ws.Cells.Item(4, 11).Value = 0.08839268234497212
$ws.Cells.Item(4, 12).Value = 0.08839268234497212
$ws.Cells.Item(4, 13).Value = 0.08839268234497212
$ws.Cells.Item(4, 14).Value = 0.08839268234497212
$ws.Cells.Item(4, 15).Value = 0.08839268234497212
$ws.Cells.Item(4, 16).Value = 0.08839268234497212
$ws.Cells.Item(4, 17).Value = 0.08839268234497212
$ws.Cells.Item(4, 18).Value = 0.08839268234497212
$ws.Cells.Item(4, 19).Value = 0.08839268234497212
$ws.Cells.Item(4, 20).Value = 0.08839268234497212
$ws.Cells.Item(4, 21).Value = 0.08839268234497212
$ws.Cells.Item(4, 22).Value = 0.08839268234497212
$ws.Cells.Item(4, 23).Value = 0.08830880999943046
$ws.Cells.Item(4, 24).Value = 0.08829813308117211
$ws.Cells.Item(4, 25).Value = 0.08829813308117211
$ws.Cells.Item(5, 3).Value = 0.5625276565551758
$ws.Cells.Item(5, 5).Value = 3983.118473920424
$ws.Cells.Item(5, 6).Value = 0.1139790135869888
$ws.Cells.Item(5, 7).Value = 0.09696133310861395
$ws.Cells.Item(5, 8).Value = 0.08851051664638104
$ws.Cells.Item(5, 9).Value = 0.08851051664638104
$ws.Cells.Item(5, 10).Value = 0.08851051664638104
$ws.Cells.Item(5, 11).Value = 0.08851051664638104
$ws.Cells.Item(5, 12).Value = 0.08851051664638104
$ws.Cells.Item(5, 13).Value = 0.08851051664638104
$ws.Cells.Item(5, 14).Value = 0.08851051664638104
$ws.Cells.Item(5, 15).Value = 0.08851051664638104
$ws.Cells.Item(5, 16).Value = 0.08851051664638104
$ws.Cells.Item(5, 17).Value = 0.08831030163587569
$ws.Cells.Item(5, 18).Value = 0.08831030163587569
$ws.Cells.Item(5, 19).Value = 0.08831030163587569
$ws.Cells.Item(5, 20).Value = 0.08831030163587569
$ws.Cells.Item(5, 21).Value = 0.08831030163587569
$ws.Cells.Item(5, 22).Value = 0.08831030163587569
$ws.Cells.Item(5, 23).Value = 0.08831030163587569
$ws.Cells.Item(5, 24).Value = 0.08831030163587569
$ws.Cells.Item(5, 25).Value = 0.08831030163587569
$ws.Cells.Item(6, 3).Value = 0.5781002044677734
$ws.Cells.Item(6, 5).Value = 3965.601166609454
$ws.Cells.Item(6, 6).Value = 0.1122850210477101
$ws.Cells.Item(6, 7).Value = 0.09386498307581029
$ws.Cells.Item(6, 8).Value = 0.08951999617818038
$ws.Cells.Item(6, 9).Value = 0.08951999617818038
$ws.Cells.Item(6, 10).Value = 0.08951999617818038
$ws.Cells.Item(6, 11).Value = 0.08951999617818038
$ws.Cells.Item(6, 12).Value = 0.08923115682709858
$ws.Cells.Item(6, 13).Value = 0.08923115682709858
$ws.Cells.Item(6, 14).Value = 0.08868089121964921
$ws.Cells.Item(6, 15).Value = 0.08868089121964921
$ws.Cells.Item(6, 16).Value = 0.08868089121964921
$ws.Cells.Item(6, 17).Value = 0.08868089121964921
$ws.Cells.Item(6, 18).Value = 0.08868089121964921
$ws.Cells.Item(6, 19).Value = 0.08868089121964921
$ws.Cells.Item(6, 20).Value = 0.08868089121964921
$ws.Cells.Item(6, 21).Value = 0.08868089121964921
$ws.Cells.Item(6, 22).Value = 0.08796883365710434
$ws.Cells.Item(6, 23).Value = 0.08796883365710434
$ws.Cells.Item(6, 24).Value = 0.08796883365710434
$ws.Cells.Item(6, 25).Value = 0.08796883365710434
$ws.Cells.Item(7, 3).Value = 0.5937752723693848
$ws.Cells.Item(7, 5).Value = 4015.866200509996
$ws.Cells.Item(7, 6).Value = 0.1024739487415298
$ws.Cells.Item(7, 7).Value = 0.09808922256392381
$ws.Cells.Item(7, 8).Value = 0.09141855512312207
$ws.Cells.Item(7, 9).Value = 0.09056566633149699
$ws.Cells.Item(7, 10).Value = 0.09056566633149699
$ws.Cells.Item(7, 11).Value = 0.09056566633149699
$ws.Cells.Item(7, 12).Value = 0.09056566633149699
$ws.Cells.Item(7, 13).Value = 0.09056566633149699
$ws.Cells.Item(7, 14).Value = 0.09031574160808997
$ws.Cells.Item(7, 15).Value = 0.09031574160808997
$ws.Cells.Item(7, 16).Value = 0.08960213956108147
$ws.Cells.Item(7, 17).Value = 0.08960213956108147
$ws.Cells.Item(7, 18).Value = 0.08960213956108147
$ws.Cells.Item(7, 19).Value = 0.08960213956108147
$ws.Cells.Item(7, 20).Value = 0.08939437953377424
$ws.Cells.Item(7, 21).Value = 0.08939437953377424
$ws.Cells.Item(7, 22).Value = 0.08922645563597494
$ws.Cells.Item(7, 23).Value = 0.08922645563597494
$ws.Cells.Item(7, 24).Value = 0.08915660016679086
$ws.Cells.Item(7, 25).Value = 0.08894865887933714
$ws.Cells.Item(8, 3).Value = 0.5781242847442627
$ws.Cells.Item(8, 5).Value = 3988.289631738367
$ws.Cells.Item(8, 6).Value = 0.1109588348104388
$ws.Cells.Item(8, 7).Value = 0.09791663166745167
$ws.Cells.Item(8, 8).Value = 0.09791663166745167
$ws.Cells.Item(8, 9).Value = 0.09244938648061855
$ws.Cells.Item(8, 10).Value = 0.09244938648061855
$ws.Cells.Item(8, 11).Value = 0.08994085561775879
$ws.Cells.Item(8, 12).Value = 0.08994085561775879
$ws.Cells.Item(8, 13).Value = 0.08994085561775879
$ws.Cells.Item(8, 14).Value = 0.08994085561775879
$ws.Cells.Item(8, 15).Value = 0.08963904298533183
$ws.Cells.Item(8, 16).Value = 0.08963904298533183
$ws.Cells.Item(8, 17).Value = 0.0885941011762193
$ws.Cells.Item(8, 18).Value = 0.0885941011762193
$ws.Cells.Item(8, 19).Value = 0.0885941011762193
$ws.Cells.Item(8, 20).Value = 0.0885941011762193
$ws.Cells.Item(8, 21).Value = 0.08841110393252175
$ws.Cells.Item(8, 22).Value = 0.08841110393252175
$ws.Cells.Item(8, 23).Value = 0.08841110393252175
$ws.Cells.Item(8, 24).Value = 0.08841110393252175
$ws.Cells.Item(8, 25).Value = 0.08841110393252175
$ws.Cells.Item(9, 3).Value = 0.5625247955322266
$ws.Cells.Item(9, 5).Value = 3992.073557720511
$ws.Cells.Item(9, 6).Value = 0.1111220854757544
$ws.Cells.Item(9, 7).Value = 0.09929990293201571
$ws.Cells.Item(9, 8).Value = 0.0918408374925941
$ws.Cells.Item(9, 9).Value = 0.08986321868277079
$ws.Cells.Item(9, 10).Value = 0.08986321868277079
$ws.Cells.Item(9, 11).Value = 0.08952076636380167
$ws.Cells.Item(9, 12).Value = 0.08952076636380167
$ws.Cells.Item(9, 13).Value = 0.08951628661716182
$ws.Cells.Item(9, 14).Value = 0.08901144645986662
$ws.Cells.Item(9, 15).Value = 0.08901144645986662
$ws.Cells.Item(9, 16).Value = 0.08901144645986662
$ws.Cells.Item(9, 17).Value = 0.08900981318984759
$ws.Cells.Item(9, 18).Value = 0.08900981318984759
$ws.Cells.Item(9, 19).Value = 0.08880996559080201
$ws.Cells.Item(9, 20).Value = 0.08880996559080201
$ws.Cells.Item(9, 21).Value = 0.08880996559080201
$ws.Cells.Item(9, 22).Value = 0.08880996559080201
$ws.Cells.Item(9, 23).Value = 0.08880996559080201
$ws.Cells.Item(9, 24).Value = 0.0887837138704558
$ws.Cells.Item(9, 25).Value = 0.08848486467291442
$ws.Cells.Item(10, 3).Value = 0.5780987739562988
$ws.Cells.Item(10, 5).Value = 3998.922586091195
$ws.Cells.Item(10, 6).Value = 0.1160189416882825
$ws.Cells.Item(10, 7).Value = 0.1067350372143929
$ws.Cells.Item(10, 8).Value = 0.09140478528807247
$ws.Cells.Item(10, 9).Value = 0.09140478528807247
$ws.Cells.Item(10, 10).Value = 0.0912114978673526
$ws.Cells.Item(10, 11).Value = 0.09069082666031124
$ws.Cells.Item(10, 12).Value = 0.08898554863375562
$ws.Cells.Item(10, 13).Value = 0.08898554863375562
$ws.Cells.Item(10, 14).Value = 0.08898554863375562
$ws.Cells.Item(10, 15).Value = 0.08898554863375562
$ws.Cells.Item(10, 16).Value = 0.08898554863375562
$ws.Cells.Item(10, 17).Value = 0.08898554863375562
$ws.Cells.Item(10, 18).Value = 0.08898554863375562
$ws.Cells.Item(10, 19).Value = 0.08898554863375562
$ws.Cells.Item(10, 20).Value = 0.08898554863375562
$ws.Cells.Item(10, 21).Value = 0.08874949782613061
$ws.Cells.Item(10, 22).Value = 0.08874949782613061
$ws.Cells.Item(10, 23).Value = 0.08874949782613061
$ws.Cells.Item(10, 24).Value = 0.08874949782613061
$ws.Cells.Item(10, 25).Value = 0.08861837399787903
$ws.Cells.Item(11, 3).Value = 0.5624749660491943
$ws.Cells.Item(11, 5).Value = 3948.226895081022
$ws.Cells.Item(11, 6).Value = 0.1139847087367801
$ws.Cells.Item(11, 7).Value = 0.101630671439467
$ws.Cells.Item(11, 8).Value = 0.08835022015274206
$ws.Cells.Item(11, 9).Value = 0.08835022015274206
$ws.Cells.Item(11, 10).Value = 0.08835022015274206
$ws.Cells.Item(11, 11).Value = 0.08835022015274206
$ws.Cells.Item(11, 12).Value = 0.08835022015274206
$ws.Cells.Item(11, 13).Value = 0.08835022015274206
$ws.Cells.Item(11, 14).Value = 0.08835022015274206
$ws.Cells.Item(11, 15).Value = 0.08835022015274206
$ws.Cells.Item(11, 16).Value = 0.08835022015274206
$ws.Cells.Item(11, 17).Value = 0.08835022015274206
$ws.Cells.Item(11, 18).Value = 0.08817625795509355
$ws.Cells.Item(11, 19).Value = 0.08817625795509355
$ws.Cells.Item(11, 20).Value = 0.08787368134451931
$ws.Cells.Item(11, 21).Value = 0.08787368134451931
$ws.Cells.Item(11, 22).Value = 0.08787368134451931
$ws.Cells.Item(11, 23).Value = 0.08765397587302197
$ws.Cells.Item(11, 24).Value = 0.08765397587302197
$ws.Cells.Item(11, 25).Value = 0.08763015390021485
